{"js": "const body = context.document.body;\n\n// Merge \"A laptop \" + \"that can connect to the \" + \"Internet\" + \" via \" into one run.\nlet results = body.search(\"A laptop that can connect to the Internet via \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"A laptop that can connect to the Internet via \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Merge \"W\" + \"i\" + \"FI\" into one run.\nresults = body.search(\"WiFI\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"WiFI\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Merge \" and ideally some familiarity \" + \"with\" + \" Python \" + \"a\" + \"n\" + \"d\" + \" \"\n// into \" and, ideally, some familiarity with Python and \".\nresults = body.search(\"and ideally some familiarity with Python and \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"and, ideally, some familiarity with Python and \", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Merge \"A laptop \" + \"that can connect to the \" + \"Internet\" + \" via \" into one run\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"A laptop that can connect to the Internet via \"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"A laptop that can connect to the Internet via \"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# Merge \"W\" + \"i\" + \"FI\" into one run\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"WiFI\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"WiFI\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# Merge \" and ideally some familiarity \" + \"with\" + \" Python \" + \"a\" + \"n\" + \"d\" + \" \"\n# into \" and, ideally, some familiarity with Python and \"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = \"and ideally some familiarity with Python and\"\n$find3.Replacement.ClearFormatting()\n$find3.Replacement.Text = \"and, ideally, some familiarity with Python and\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2)\n"}
